$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "A"
$ws.Range("A2").Value = "b"
$ws.Range("A3:A12").Value = "c"

$ws.Range("B12").Select()
